# Update points 09876543 -> 0.00
#
# A new customer row (phone entered as a plain number, 9876543 - losing the
# leading zero that the existing "09876543" text entry had) is inserted
# above the existing row for "09876543" on Sheet1. This pushes the old
# "09876543" row down by one, from row 22 to row 23, and adds a fresh
# row 22 with phone = 9876543 (numeric) and total_points = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "09876543" row (and everything at/after it) down one
# row, opening up a blank row 22 for the new entry.
$ws.Rows("22:22").Insert()

# New row 22: phone captured as a genuine number (no leading zero),
# blank birthday, 0 points - mirrors the other freshly-added 9876543 rows
# above it.
$ws.Range("A22").Value = 9876543
# A lone "'" forces the cell into text mode with an empty value (matching
# the blank-but-text birthday cells used throughout this column), then
# resetting to the "Normal" style drops the quote-prefix formatting that
# entering a leading apostrophe would otherwise leave behind.
$ws.Range("B22").Value = "'"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = 0
